# EventController Coverage Change Role Implemetation ModifyEventDetails not working
#
# This script updates two "modifyUserProfile" test-case rows (85 and 86) on
# Sheet1 of the UserTestCases workbook so that the username/first/last name
# and address fields are left blank (simulating a "modify profile with
# blank fields" scenario) and the corresponding error-message columns are
# populated, and also refreshes the sheet's active selection/scroll.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 85 (Test Case 83) - role "User"
# ---------------------------------------------------------------------
# username (C85) and password (D85) stay as-is; first/last name are cleared
$ws.Range("E85").Value = ""
$ws.Range("F85").Value = ""
# phone number gets shortened by one digit
$ws.Range("I85").Value = 469258056
# email (J85) is cleared but keeps its existing style
$ws.Range("J85").Value = ""
# street number / street name / city / state cleared
$ws.Range("K85").Value = ""
$ws.Range("L85").Value = ""
$ws.Range("M85").Value = ""
$ws.Range("N85").Value = ""
# zipcode shortened
$ws.Range("O85").Value = 7601
# error columns populated
$ws.Range("P85").Value = "Please correct the following errors"
$ws.Range("S85").Value = "First Name can not be blank."
$ws.Range("T85").Value = "Last Name can not be blank."
$ws.Range("W85").Value = "Phone number must have 10 digits"
$ws.Range("X85").Value = "Email address cannot be empty"
$ws.Range("Y85").Value = "Street number length must be >0 and <7."
$ws.Range("Z85").Value = "Street name length must be >0 and <40."
$ws.Range("AA85").Value = "City Name cannot be Empty"
$ws.Range("AB85").Value = "State Cannot be Empty."
$ws.Range("AC85").Value = "Zipcode must have a length of 5"

# ---------------------------------------------------------------------
# Row 86 (Test Case 84) - role "Admin"
# ---------------------------------------------------------------------
# username (C86) is cleared entirely, password (D86) stays
$ws.Range("C86").Value = ""
$ws.Range("E86").Value = ""
$ws.Range("F86").Value = ""
# role changes from User to Admin
$ws.Range("G86").Value = "Admin"
# phone number gets shortened by one digit
$ws.Range("I86").Value = 469258056
# email (J86) is cleared but keeps its existing style
$ws.Range("J86").Value = ""
# street number / street name / city / state cleared
$ws.Range("K86").Value = ""
$ws.Range("L86").Value = ""
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = ""
# zipcode shortened
$ws.Range("O86").Value = 7601
# error columns populated
$ws.Range("P86").Value = "Please correct the following errors"
$ws.Range("Q86").Value = "Username can not be blank."
$ws.Range("S86").Value = "First Name can not be blank."
$ws.Range("T86").Value = "Last Name can not be blank."
$ws.Range("U86").Value = "There can only be one admin"
$ws.Range("W86").Value = "Phone number must have 10 digits"
$ws.Range("X86").Value = "Email address cannot be empty"
$ws.Range("Y86").Value = "Street number length must be >0 and <7."
$ws.Range("Z86").Value = "Street name length must be >0 and <40."
$ws.Range("AA86").Value = "City Name cannot be Empty"
$ws.Range("AB86").Value = "State Cannot be Empty."
$ws.Range("AC86").Value = "Zipcode must have a length of 5"

# ---------------------------------------------------------------------
# Sheet view: scroll so row 66 / column Q is the top-left visible cell,
# and select S86:AC86 (active cell S86) to match where the edits were made.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 66
$win.ScrollColumn = 17
$ws.Range("S86:AC86").Select()
